# Updated cryptos list (price + 1h volume change %) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.919.95"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "3.415.32"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.91"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.56"
$ws.Range("E6").Value = "  -4.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +4.84%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +11.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("E10").Value = "  +15.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.83"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.141"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.22"
$ws.Range("E13").Value = "  +7.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.83"
$ws.Range("E14").Value = "  +4.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000202"
$ws.Range("E15").Value = "  +58.86%  "

$ws.Range("D16").Value = "3.429.74"
$ws.Range("E16").Value = "  +2.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.68"
$ws.Range("E17").Value = "  +15.22%  "

$ws.Range("E18").Value = "  +3.97%  "

$ws.Range("D19").Value = "61.981.63"
$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "404.88"
$ws.Range("E20").Value = "  +28.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "90.72"
$ws.Range("E21").Value = "  +6.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.18"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.39"
$ws.Range("E23").Value = "  +4.47%  "

$ws.Range("E24").Value = "  +2.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "33.00"
$ws.Range("E25").Value = "  +11.77%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.50"
$ws.Range("E27").Value = "  +1.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.61"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.70"
$ws.Range("E29").Value = "  +5.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "43.88"
$ws.Range("E32").Value = "  +7.73%  "

$ws.Range("E33").Value = "  +3.87%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0498"
$ws.Range("E35").Value = "  +3.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.85"
$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.39"
$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("E39").Value = "  -1.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.132"
$ws.Range("E40").Value = "  +6.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.313"
$ws.Range("E41").Value = "  +6.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.82"
$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.00"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("E45").Value = "  +6.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.74"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.93"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("D48").Value = "2.106.24"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("E49").Value = "  +8.45%  "

$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.126"
$ws.Range("E51").Value = "  +12.64%  "
